$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.502.83"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "2.221.94"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.11%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.12%  "
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "2.557.86"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "2.219.57"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "42.437.41"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.24%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "237.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.17%  "
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.01%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  -6.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0877"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.74%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.126"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.76%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.40%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("E41").Value = "  -5.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.230"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.59%  "
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("E47").Value = "  -6.10%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
